# "NC files from RMI"
# The "Cost Data" sheet's annual-total rows (B88/C88 and B96) were being
# reported on a per-decade (or similarly scaled) basis; this divides them
# by 10 to correct the units, and applies a 2-decimal dollar number format
# to the corrected cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cost Data")

# Scale the totals down by a factor of 10.
$ws.Range("B88").Formula = "=B54/10"
$ws.Range("C88").Formula = "=B55/10"
$ws.Range("B96").Formula = "=B87/10"

# Switch the corrected cells to a "$#,##0.00" number format (was "$#,##0").
$ws.Range("B88").NumberFormat = """$""#,##0.00"
$ws.Range("C88").NumberFormat = """$""#,##0.00"
$ws.Range("B96").NumberFormat = """$""#,##0.00"
